# Refactor code, implement activities, begin unit tests
#
# Updates the "Integré" (D column) status from "Non" to "Oui" for several
# feature rows, updates the "Interet" (B column) priority for row 25 from
# "Faible" to "Nul", and moves the active cell selection to D8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Integré" column (D) flips from "Non" to "Oui"
$rowsNowIntegrated = @(9, 12, 13, 14, 15, 16, 17, 24)
foreach ($row in $rowsNowIntegrated) {
    $ws.Range("D$row").Value = "Oui"
}

# Row 25's "Interet" (B) priority changes from "Faible" to "Nul"
$ws.Range("B25").Value = "Nul"

# Move the current selection to D8
[void]$ws.Range("D8").Select()
